$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.829.18"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "3.789.80"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.85"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.60"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").Value = "3.785.05"
$ws.Range("E7").Value = "  -2.18%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +5.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.88"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "4.433.28"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "3.799.61"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.12"
$ws.Range("E17").Value = "  +5.45%  "
$ws.Range("D18").Value = "67.946.59"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.27"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.58"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.56"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000150"
$ws.Range("E24").Value = "  -5.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.48"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.22"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.34"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.93"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").Value = "3.945.94"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.64"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.49"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.19"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").Value = "3.758.29"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.105"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.77"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.91"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.319"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.75"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "407.48"
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.25"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000277"
$ws.Range("E49").Value = "  -7.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.98"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0356"
$ws.Range("E51").Value = "  +0.14%  "
